$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.072117533251326
$ws.Cells.Item(2, 4).Value = 1.056029215550515
$ws.Cells.Item(2, 5).Value = 1.073147793998621
$ws.Cells.Item(2, 6).Value = 1.076381399202439
$ws.Cells.Item(2, 9).Value = 1.039395166223192
$ws.Cells.Item(2, 10).Value = 1.0770386973822
$ws.Cells.Item(2, 11).Value = 1.058768049102491
$ws.Cells.Item(2, 12).Value = 1.075840424706898
$ws.Cells.Item(2, 13).Value = 1.079065480798364
$ws.Cells.Item(2, 14).Value = 1.078568215870549

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.074568496578414
$ws.Cells.Item(3, 4).Value = 1.05708843796575
$ws.Cells.Item(3, 5).Value = 1.075190035603059
$ws.Cells.Item(3, 6).Value = 1.078156078667016
$ws.Cells.Item(3, 9).Value = 1.039662431027937
$ws.Cells.Item(3, 10).Value = 1.079140837055604
$ws.Cells.Item(3, 11).Value = 1.059640203878783
$ws.Cells.Item(3, 12).Value = 1.077696533195291
$ws.Cells.Item(3, 13).Value = 1.080655314720308
$ws.Cells.Item(3, 14).Value = 1.080673340823407

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.076148308348082
$ws.Cells.Item(4, 4).Value = 1.057770180046185
$ws.Cells.Item(4, 5).Value = 1.07650596298332
$ws.Cells.Item(4, 6).Value = 1.079299187053032
$ws.Cells.Item(4, 9).Value = 1.039832611500676
$ws.Cells.Item(4, 10).Value = 1.080494902998806
$ws.Cells.Item(4, 11).Value = 1.060200308812846
$ws.Cells.Item(4, 12).Value = 1.078891630898569
$ws.Cells.Item(4, 13).Value = 1.081678368314653
$ws.Cells.Item(4, 14).Value = 1.082029329695562

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.076811034637878
$ws.Cells.Item(5, 4).Value = 1.058055924334203
$ws.Cells.Item(5, 5).Value = 1.077057884259242
$ws.Cells.Item(5, 6).Value = 1.079778524011968
$ws.Cells.Item(5, 9).Value = 1.039903499780003
$ws.Cells.Item(5, 10).Value = 1.081062711958868
$ws.Cells.Item(5, 11).Value = 1.060434774763707
$ws.Cells.Item(5, 12).Value = 1.079392660489513
$ws.Cells.Item(5, 13).Value = 1.082107125802797
$ws.Cells.Item(5, 14).Value = 1.082597945009476

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.07692222679356
$ws.Cells.Item(6, 4).Value = 1.05810385198803
$ws.Cells.Item(6, 5).Value = 1.077150479265982
$ws.Cells.Item(6, 6).Value = 1.079858935658436
$ws.Cells.Item(6, 9).Value = 1.039915363922877
$ws.Cells.Item(6, 10).Value = 1.081157966103096
$ws.Cells.Item(6, 11).Value = 1.060474084216312
$ws.Cells.Item(6, 12).Value = 1.079476704943884
$ws.Cells.Item(6, 13).Value = 1.082179038536146
$ws.Cells.Item(6, 14).Value = 1.082693334425514

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.076157169281292
$ws.Cells.Item(7, 4).Value = 1.057774001537697
$ws.Cells.Item(7, 5).Value = 1.076513342826926
$ws.Cells.Item(7, 6).Value = 1.079305596758874
$ws.Cells.Item(7, 9).Value = 1.039833561283096
$ws.Cells.Item(7, 10).Value = 1.080502495703821
$ws.Cells.Item(7, 11).Value = 1.060203445679538
$ws.Cells.Item(7, 12).Value = 1.078898331100169
$ws.Cells.Item(7, 13).Value = 1.081684102603015
$ws.Cells.Item(7, 14).Value = 1.082036933183089

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.072947140242529
$ws.Cells.Item(8, 4).Value = 1.056387947282798
$ws.Cells.Item(8, 5).Value = 1.073839145238264
$ws.Cells.Item(8, 6).Value = 1.076982258525851
$ws.Cells.Item(8, 9).Value = 1.039486063421911
$ws.Cells.Item(8, 10).Value = 1.077750421264745
$ws.Cells.Item(8, 11).Value = 1.059063682767977
$ws.Cells.Item(8, 12).Value = 1.076468950903805
$ws.Cells.Item(8, 13).Value = 1.07960396252623
$ws.Cells.Item(8, 14).Value = 1.079280950482642

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.067241938691486
$ws.Cells.Item(9, 4).Value = 1.053917044046148
$ws.Cells.Item(9, 5).Value = 1.069083009617572
$ws.Cells.Item(9, 6).Value = 1.072847031352683
$ws.Cells.Item(9, 9).Value = 1.038852396793692
$ws.Cells.Item(9, 10).Value = 1.072852196905846
$ws.Cells.Item(9, 11).Value = 1.05702226091976
$ws.Cells.Item(9, 12).Value = 1.072141313265671
$ws.Cells.Item(9, 13).Value = 1.075893915709763
$ws.Cells.Item(9, 14).Value = 1.07437577008332

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.063403164169072
$ws.Cells.Item(10, 4).Value = 1.052249796295514
$ws.Cells.Item(10, 5).Value = 1.065880743629884
$ws.Cells.Item(10, 6).Value = 1.070060855631946
$ws.Cells.Item(10, 9).Value = 1.038415312289912
$ws.Cells.Item(10, 10).Value = 1.069551775185387
$ws.Cells.Item(10, 11).Value = 1.055638323530943
$ws.Cells.Item(10, 12).Value = 1.069222872308254
$ws.Cells.Item(10, 13).Value = 1.073389015603093
$ws.Cells.Item(10, 14).Value = 1.071070661385454

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.061731983844769
$ws.Cells.Item(11, 4).Value = 1.05152293558188
$ws.Cells.Item(11, 5).Value = 1.064486198167048
$ws.Cells.Item(11, 6).Value = 1.068847075100166
$ws.Cells.Item(11, 9).Value = 1.03822250969699
$ws.Cells.Item(11, 10).Value = 1.068113877120114
$ws.Cells.Item(11, 11).Value = 1.055033431113018
$ws.Cells.Item(11, 12).Value = 1.067950814532179
$ws.Cells.Item(11, 13).Value = 1.072296533029806
$ws.Cells.Item(11, 14).Value = 1.069630721339998

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.061109836201102
$ws.Cells.Item(12, 4).Value = 1.051252190051915
$ws.Cells.Item(12, 5).Value = 1.063966968931906
$ws.Cells.Item(12, 6).Value = 1.068395086588979
$ws.Cells.Item(12, 9).Value = 1.038150356111141
$ws.Cells.Item(12, 10).Value = 1.067578413843836
$ws.Cells.Item(12, 11).Value = 1.054807884339463
$ws.Cells.Item(12, 12).Value = 1.067477023157876
$ws.Cells.Item(12, 13).Value = 1.071889527075236
$ws.Cells.Item(12, 14).Value = 1.06909449764445

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.061243353005482
$ws.Cells.Item(13, 4).Value = 1.051310300411486
$ws.Cells.Item(13, 5).Value = 1.064078401803798
$ws.Cells.Item(13, 6).Value = 1.068492091610193
$ws.Cells.Item(13, 9).Value = 1.038165857748304
$ws.Cells.Item(13, 10).Value = 1.067693334950286
$ws.Cells.Item(13, 11).Value = 1.054856304187061
$ws.Cells.Item(13, 12).Value = 1.067578712112935
$ws.Cells.Item(13, 13).Value = 1.071976886457706
$ws.Cells.Item(13, 14).Value = 1.069209581952053

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.061680585714042
$ws.Cells.Item(14, 4).Value = 1.051500571211429
$ws.Cells.Item(14, 5).Value = 1.064443303900238
$ws.Cells.Item(14, 6).Value = 1.068809736960398
$ws.Cells.Item(14, 9).Value = 1.038216556471434
$ws.Cells.Item(14, 10).Value = 1.068069643632416
$ws.Cells.Item(14, 11).Value = 1.055014805039087
$ws.Cells.Item(14, 12).Value = 1.06791167737478
$ws.Cells.Item(14, 13).Value = 1.072262914639824
$ws.Cells.Item(14, 14).Value = 1.069586425035674

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.061949792610731
$ws.Cells.Item(15, 4).Value = 1.051617702518011
$ws.Cells.Item(15, 5).Value = 1.064667967549525
$ws.Cells.Item(15, 6).Value = 1.069005297183187
$ws.Cells.Item(15, 9).Value = 1.038247722141917
$ws.Cells.Item(15, 10).Value = 1.068301317806929
$ws.Cells.Item(15, 11).Value = 1.055112347841617
$ws.Cells.Item(15, 12).Value = 1.068116655863966
$ws.Cells.Item(15, 13).Value = 1.072438984870876
$ws.Cells.Item(15, 14).Value = 1.069818428214089

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.063513881771619
$ws.Cells.Item(16, 4).Value = 1.052297930371714
$ws.Cells.Item(16, 5).Value = 1.065973124303844
$ws.Cells.Item(16, 6).Value = 1.070141252680189
$ws.Cells.Item(16, 9).Value = 1.038428032829143
$ws.Cells.Item(16, 10).Value = 1.069647014839828
$ws.Cells.Item(16, 11).Value = 1.055678348121745
$ws.Cells.Item(16, 12).Value = 1.069307115412654
$ws.Cells.Item(16, 13).Value = 1.073461352190204
$ws.Cells.Item(16, 14).Value = 1.071166036291127

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.064492558305266
$ws.Cells.Item(17, 4).Value = 1.052723287708454
$ws.Cells.Item(17, 5).Value = 1.066789660326805
$ws.Cells.Item(17, 6).Value = 1.070851817573709
$ws.Cells.Item(17, 9).Value = 1.038540184467893
$ws.Cells.Item(17, 10).Value = 1.07048875077871
$ws.Cells.Item(17, 11).Value = 1.056031864604597
$ws.Cells.Item(17, 12).Value = 1.070051596527059
$ws.Cells.Item(17, 13).Value = 1.074100534518731
$ws.Cells.Item(17, 14).Value = 1.072008967591592

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.06506254213088
$ws.Cells.Item(18, 4).Value = 1.052970916745681
$ws.Cells.Item(18, 5).Value = 1.067265168284695
$ws.Cells.Item(18, 6).Value = 1.071265571292759
$ws.Cells.Item(18, 9).Value = 1.038605259266844
$ws.Cells.Item(18, 10).Value = 1.070978875575351
$ws.Cells.Item(18, 11).Value = 1.056237521840286
$ws.Cells.Item(18, 12).Value = 1.07048503632943
$ws.Cells.Item(18, 13).Value = 1.074472604002821
$ws.Cells.Item(18, 14).Value = 1.072499788421662

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.065256747232748
$ws.Cells.Item(19, 4).Value = 1.053055271815959
$ws.Cells.Item(19, 5).Value = 1.067427175862316
$ws.Cells.Item(19, 6).Value = 1.071406531711492
$ws.Cells.Item(19, 9).Value = 1.038627390370253
$ws.Cells.Item(19, 10).Value = 1.071145853329921
$ws.Cells.Item(19, 11).Value = 1.056307554019724
$ws.Cells.Item(19, 12).Value = 1.070632692947183
$ws.Cells.Item(19, 13).Value = 1.074599342984809
$ws.Cells.Item(19, 14).Value = 1.072667003303796

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.064387645056851
$ws.Cells.Item(20, 4).Value = 1.05267770010872
$ws.Cells.Item(20, 5).Value = 1.066702133015527
$ws.Cells.Item(20, 6).Value = 1.070775654021884
$ws.Cells.Item(20, 9).Value = 1.038528187009121
$ws.Cells.Item(20, 10).Value = 1.070398528243408
$ws.Cells.Item(20, 11).Value = 1.055993991903986
$ws.Cells.Item(20, 12).Value = 1.069971804121314
$ws.Cells.Item(20, 13).Value = 1.074032034592803
$ws.Cells.Item(20, 14).Value = 1.071918616929941

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.061551870544376
$ws.Cells.Item(21, 4).Value = 1.051444562234516
$ws.Cells.Item(21, 5).Value = 1.064335883684934
$ws.Cells.Item(21, 6).Value = 1.068716229954715
$ws.Cells.Item(21, 9).Value = 1.038201641864216
$ws.Cells.Item(21, 10).Value = 1.067958868102683
$ws.Cells.Item(21, 11).Value = 1.054968154433618
$ws.Cells.Item(21, 12).Value = 1.067813663386916
$ws.Cells.Item(21, 13).Value = 1.072178720101295
$ws.Cells.Item(21, 14).Value = 1.069475492191983

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.059760799375848
$ws.Cells.Item(22, 4).Value = 1.050664850920086
$ws.Cells.Item(22, 5).Value = 1.062840975165499
$ws.Cells.Item(22, 6).Value = 1.067414797082545
$ws.Cells.Item(22, 9).Value = 1.037993213547696
$ws.Cells.Item(22, 10).Value = 1.066417044016402
$ws.Cells.Item(22, 11).Value = 1.054318169831147
$ws.Cells.Item(22, 12).Value = 1.066449257725875
$ws.Cells.Item(22, 13).Value = 1.071006455524465
$ws.Cells.Item(22, 14).Value = 1.06793147853865

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.060711065451793
$ws.Cells.Item(23, 4).Value = 1.051078612122228
$ws.Cells.Item(23, 5).Value = 1.063634145947922
$ws.Cells.Item(23, 6).Value = 1.068105347024002
$ws.Cells.Item(23, 9).Value = 1.03810400276264
$ws.Cells.Item(23, 10).Value = 1.067235158674532
$ws.Cells.Item(23, 11).Value = 1.054663218216248
$ws.Cells.Item(23, 12).Value = 1.067173278287749
$ws.Cells.Item(23, 13).Value = 1.071628570254778
$ws.Cells.Item(23, 14).Value = 1.068750755013434

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.064435053498163
$ws.Cells.Item(24, 4).Value = 1.052698300649076
$ws.Cells.Item(24, 5).Value = 1.066741685198081
$ws.Cells.Item(24, 6).Value = 1.070810071235163
$ws.Cells.Item(24, 9).Value = 1.038533609198901
$ws.Cells.Item(24, 10).Value = 1.070439298534404
$ws.Cells.Item(24, 11).Value = 1.056011106624348
$ws.Cells.Item(24, 12).Value = 1.07000786135589
$ws.Cells.Item(24, 13).Value = 1.074062989091749
$ws.Cells.Item(24, 14).Value = 1.071959445119426

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.068722903474183
$ws.Cells.Item(25, 4).Value = 1.054559292245091
$ws.Cells.Item(25, 5).Value = 1.070317988007733
$ws.Cells.Item(25, 6).Value = 1.073921138330306
$ws.Cells.Item(25, 9).Value = 1.039018772984455
$ws.Cells.Item(25, 10).Value = 1.074124501097587
$ws.Cells.Item(25, 11).Value = 1.057554010415
$ws.Cells.Item(25, 12).Value = 1.073265847228013
$ws.Cells.Item(25, 13).Value = 1.076858493566837
$ws.Cells.Item(25, 14).Value = 1.075649881092947
